$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 21 de Abril de 2020 a las 01:22"

# Re-rank countries: Jamaica, Gabon and Sudan moved up in the standings,
# pushing the countries below them down by one row each.
$ws.Range("A125").Value = "Jamaica"
$ws.Range("A126").Value = "El Salvador"
$ws.Range("A127").Value = "Paraguay"
$ws.Range("A137").Value = "Gabon"
$ws.Range("A138").Value = "Birmania"
$ws.Range("A139").Value = "Trinidad yTobago"
$ws.Range("A140").Value = "Etiopia"
$ws.Range("A141").Value = "Sudan"
$ws.Range("A142").Value = "Liberia"
$ws.Range("A143").Value = "Aruba"
$ws.Range("A144").Value = "Guayana Francesa"
$ws.Range("A145").Value = "Monaco"

# Updated case counts for the affected rows
$ws.Range("B4").Value = 791625
$ws.Range("C4").Value = 26989
$ws.Range("D4").Value = 71895
$ws.Range("E4").Value = 677272
$ws.Range("F4").Value = 13908
$ws.Range("G4").Value = 1883
$ws.Range("H4").Value = 42458
$ws.Range("B8").Value = 147065
$ws.Range("C8").Value = 1323
$ws.Range("D8").Value = 91500
$ws.Range("E8").Value = 50703
$ws.Range("F8").Value = 2889
$ws.Range("G8").Value = 220
$ws.Range("H8").Value = 4862
$ws.Range("B14").Value = 40743
$ws.Range("C14").Value = 2089
$ws.Range("D14").Value = 22130
$ws.Range("E14").Value = 16026
$ws.Range("F14").Value = 7919
$ws.Range("G14").Value = 125
$ws.Range("H14").Value = 2587
$ws.Range("B125").Value = 223
$ws.Range("C125").Value = 50
$ws.Range("D125").Value = 27
$ws.Range("E125").Value = 191
$ws.Range("F125").Value = 0
$ws.Range("G125").Value = 0
$ws.Range("H125").Value = 5
$ws.Range("B126").Value = 218
$ws.Range("C126").Value = 17
$ws.Range("D126").Value = 46
$ws.Range("E126").Value = 165
$ws.Range("F126").Value = 2
$ws.Range("G126").Value = 0
$ws.Range("H126").Value = 7
$ws.Range("B127").Value = 208
$ws.Range("C127").Value = 2
$ws.Range("D127").Value = 46
$ws.Range("E127").Value = 154
$ws.Range("F127").Value = 1
$ws.Range("G127").Value = 0
$ws.Range("H127").Value = 8
$ws.Range("B137").Value = 120
$ws.Range("C137").Value = 11
$ws.Range("D137").Value = 7
$ws.Range("E137").Value = 112
$ws.Range("F137").Value = 2
$ws.Range("G137").Value = 0
$ws.Range("H137").Value = 1
$ws.Range("B138").Value = 117
$ws.Range("C138").Value = 6
$ws.Range("D138").Value = 7
$ws.Range("E138").Value = 105
$ws.Range("F138").Value = 0
$ws.Range("G138").Value = 0
$ws.Range("H138").Value = 5
$ws.Range("B139").Value = 114
$ws.Range("C139").Value = 0
$ws.Range("D139").Value = 21
$ws.Range("E139").Value = 85
$ws.Range("F139").Value = 0
$ws.Range("G139").Value = 0
$ws.Range("H139").Value = 8
$ws.Range("B140").Value = 111
$ws.Range("C140").Value = 3
$ws.Range("D140").Value = 16
$ws.Range("E140").Value = 92
$ws.Range("F140").Value = 1
$ws.Range("G140").Value = 0
$ws.Range("H140").Value = 3
$ws.Range("B141").Value = 107
$ws.Range("C141").Value = 41
$ws.Range("D141").Value = 8
$ws.Range("E141").Value = 87
$ws.Range("F141").Value = 0
$ws.Range("G141").Value = 2
$ws.Range("H141").Value = 12
$ws.Range("B142").Value = 99
$ws.Range("C142").Value = 8
$ws.Range("D142").Value = 7
$ws.Range("E142").Value = 84
$ws.Range("F142").Value = 0
$ws.Range("G142").Value = 0
$ws.Range("H142").Value = 8
$ws.Range("B143").Value = 97
$ws.Range("C143").Value = 0
$ws.Range("D143").Value = 49
$ws.Range("E143").Value = 46
$ws.Range("F143").Value = 4
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 2
$ws.Range("B144").Value = 97
$ws.Range("C144").Value = 1
$ws.Range("D144").Value = 69
$ws.Range("E144").Value = 27
$ws.Range("F144").Value = 2
$ws.Range("G144").Value = 1
$ws.Range("H144").Value = 1
$ws.Range("B145").Value = 94
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 23
$ws.Range("E145").Value = 68
$ws.Range("F145").Value = 3
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 3
$ws.Range("D161").Value = 15
$ws.Range("E161").Value = 35
